$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values per diff
$ws.Range("D7").Value = 95
$ws.Range("G7").Value = 98
$ws.Range("G22").Value = 99
$ws.Range("H22").Value = 140

# Update selected cell / active cell in the view
$ws.Range("G8").Select()
